$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")

# Row 12: "Functions"/"USDLaunch" -> "USD"/"Launch"
$ws.Range("C12").Value = "USD"
$ws.Range("D12").Value = "Launch"

# Row 14: "Functions"/"CrmChangeArea" -> "Crm"/"ChangeArea"
$ws.Range("C14").Value = "Crm"
$ws.Range("D14").Value = "ChangeArea"

# Row 15: "Functions"/"CrmOpenEntity" -> "Crm"/"OpenEntity"
$ws.Range("C15").Value = "Crm"
$ws.Range("D15").Value = "OpenEntity"

# Row 16: "Functions"/"CrmClickButton" -> "Crm"/"ClickButton"
$ws.Range("C16").Value = "Crm"
$ws.Range("D16").Value = "ClickButton"

# Row 20: "Functions"/"CrmLookupField" -> "Crm"/"LookupField"
$ws.Range("C20").Value = "Crm"
$ws.Range("D20").Value = "LookupField"

# Row 32: "Functions"/"USDClose" -> "USD"/"Close"
$ws.Range("C32").Value = "USD"
$ws.Range("D32").Value = "Close"
